$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl_testcases")

# --- Update TCID values (column B) -- the formulas in column A recompute automatically ---
# Write B3 before B2 so the new shared-strings end up in the same order as the target file
# (test_tc2 then test_tc1).
$ws.Range("B3").Value = "test_tc2"
$ws.Range("B2").Value = "test_tc1"

# --- Merge the three separate "containsText" FLow rules (on C1, D1:F1 and A1:B1) into a
# single rule covering A1:F1, reusing the A1:B1 rule object so its formatting/dxf survive ---
$fcFlow = $ws.Range("A1:B1").FormatConditions.Item(1)
$fcFlow.ModifyAppliesToRange($ws.Range("A1:F1"))
$ws.Range("C1").FormatConditions.Delete()
$ws.Range("D1:F1").FormatConditions.Delete()

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()

$wb.Save()
